$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns keep their exact text representation
# (many price strings look numeric, e.g. '1.010', '0.000008786', and would
# otherwise be silently reinterpreted/rounded as floating point numbers).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.182.24'
$ws.Range("E2").Value = '  +0.91%  '

# Row 3
$ws.Range("D3").Value = '1.833.35'
$ws.Range("E3").Value = '  +0.85%  '

# Row 4
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +1.08%  '

# Row 5
$ws.Range("D5").Value = '313.39'
$ws.Range("E5").Value = '  +1.04%  '

# Row 6
$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.93%  '

# Row 7
$ws.Range("D7").Value = '0.4705'
$ws.Range("E7").Value = '  +0.58%  '

# Row 8
$ws.Range("D8").Value = '0.3686'
$ws.Range("E8").Value = '  -0.52%  '

# Row 9
$ws.Range("E9").Value = '  +0.43%  '

# Row 10
$ws.Range("D10").Value = '0.8815'
$ws.Range("E10").Value = '  +1.15%  '

# Row 11
$ws.Range("D11").Value = '20.44'
$ws.Range("E11").Value = '  -0.03%  '

# Row 12
$ws.Range("D12").Value = '1.842.30'
$ws.Range("E12").Value = '  +0.09%  '

# Row 13
$ws.Range("D13").Value = '0.07336'
$ws.Range("E13").Value = '  +3.68%  '

# Row 14
$ws.Range("D14").Value = '5.478'
$ws.Range("E14").Value = '  +1.98%  '

# Row 15
$ws.Range("D15").Value = '92.83'
$ws.Range("E15").Value = '  +0.75%  '

# Row 16
$ws.Range("D16").Value = '6.563'
$ws.Range("E16").Value = '  +0.73%  '

# Row 17
$ws.Range("E17").Value = '  +1.05%  '

# Row 18
$ws.Range("D18").Value = '0.000008786'
$ws.Range("E18").Value = '  +0.62%  '

# Row 19
$ws.Range("D19").Value = '1.009'
$ws.Range("E19").Value = '  +0.88%  '

# Row 20
$ws.Range("D20").Value = '14.79'
$ws.Range("E20").Value = '  +0.29%  '

# Row 21
$ws.Range("D21").Value = '27.207.21'
$ws.Range("E21").Value = '  +0.87%  '

# Row 22
$ws.Range("D22").Value = '5.308'
$ws.Range("E22").Value = '  -0.90%  '

# Row 23
$ws.Range("E23").Value = '  +0.98%  '

# Row 24
$ws.Range("D24").Value = '2.065.19'
$ws.Range("E24").Value = '  -0.45%  '

# Row 25
$ws.Range("D25").Value = '1.902'
$ws.Range("E25").Value = '  +0.26%  '

# Row 26
$ws.Range("D26").Value = '152.48'

# Row 27
$ws.Range("D27").Value = '18.53'
$ws.Range("E27").Value = '  +0.88%  '

# Row 28
$ws.Range("D28").Value = '2.157'
$ws.Range("E28").Value = '  -2.02%  '

# Row 29
$ws.Range("D29").Value = '5.271'
$ws.Range("E29").Value = '  -0.78%  '

# Row 30
$ws.Range("D30").Value = '117.67'
$ws.Range("E30").Value = '  +1.88%  '

# Row 31
$ws.Range("D31").Value = '0.08927'
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("D32").Value = '0.7604'
$ws.Range("E32").Value = '  -0.84%  '

# Row 33
$ws.Range("D33").Value = '1.171'
$ws.Range("E33").Value = '  +0.66%  '

# Row 34
$ws.Range("D34").Value = '4.545'
$ws.Range("E34").Value = '  +1.27%  '

# Row 35
$ws.Range("E35").Value = '  +0.59%  '

# Row 36
$ws.Range("E36").Value = '  +0.99%  '

# Row 37
$ws.Range("D37").Value = '1.102'
$ws.Range("E37").Value = '  +0.26%  '

# Row 38
$ws.Range("D38").Value = '0.05346'
$ws.Range("E38").Value = '  +1.50%  '

# Row 39
$ws.Range("D39").Value = '0.01960'
$ws.Range("E39").Value = '  -0.15%  '

# Row 40
$ws.Range("D40").Value = '3.002'
$ws.Range("E40").Value = '  +2.06%  '

# Row 41
$ws.Range("D41").Value = '7.338'
$ws.Range("E41").Value = '  +0.76%  '

# Row 42
$ws.Range("D42").Value = '2.411'
$ws.Range("E42").Value = '  +1.73%  '

# Row 43
$ws.Range("D43").Value = '0.5342'
$ws.Range("E43").Value = '  -0.95%  '

# Row 44
$ws.Range("D44").Value = '0.1663'
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("D45").Value = '8.538'
$ws.Range("E45").Value = '  +0.58%  '

# Row 46
$ws.Range("D46").Value = '0.4942'
$ws.Range("E46").Value = '  -0.59%  '

# Row 47
$ws.Range("D47").Value = '10.57'
$ws.Range("E47").Value = '  +1.56%  '

# Row 48
$ws.Range("D48").Value = '1.010'
$ws.Range("E48").Value = '  +0.99%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '103.91'
$ws.Range("E49").Value = '  +0.71%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.670'
$ws.Range("E50").Value = '  -0.51%  '

# Row 51
$ws.Range("D51").Value = '0.06328'
$ws.Range("E51").Value = '  +0.67%  '
